$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "reviews_count" column (column E) is removed entirely; all columns
# to its right (F:K -> reviews_average, latitude, longitude,
# is_permanently_closed, gmaps_link, latest_review_date) shift one
# position to the left (E:J).
$ws.Columns.Item(5).Delete()
